$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Replace "Ready for handoff" text with "In Translation" wherever it occurs
$overview.Range("E2:F3").Replace("Ready for handoff", "In Translation")
$zhcn.Range("C2:C3").Replace("Ready for handoff", "In Translation")
$dede.Range("C2:C3").Replace("Ready for handoff", "In Translation")

# Adjust column widths
$overview.Range("E1").EntireColumn.ColumnWidth = 13.4101845877511
$overview.Range("F1").EntireColumn.ColumnWidth = 13.4101845877511
$zhcn.Range("C1").EntireColumn.ColumnWidth = 13.4101845877511
$dede.Range("C1").EntireColumn.ColumnWidth = 13.4101845877511
